$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Grønne bananer" -> "Grønne oliven"
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Grønne bananer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Grønne oliven", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the end of the document (after
#    "... En tannbørste ") to wrap the "Små grønne bananer" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -eq "Små grønne bananer" + [char]13) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $d.Bookmarks.Add("_GoBack", $target.Range) | Out-Null
}
